# Roughing out the yaml format for data definitions
#
# - Adds a "Customer_Location" header to the Customers sheet (new column C)
# - Makes "Locations" the active sheet/tab (was "Sales")

$wb = $excel.ActiveWorkbook

$customers = $wb.Worksheets.Item("Customers")
$customers.Range("C1").Value = "Customer_Location"
$customers.Range("C1").Select()

$locations = $wb.Worksheets.Item("Locations")
$locations.Activate()
